# Update the "want to go" counts (column F) for two events that appear
# both on the "展览" sheet and on the aggregated "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F6").Value = 26
    $ws.Range("F9").Value = 235
}
